$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "'10-JAN-26"
$ws.Range("B2").Value = "SM-454"
$ws.Range("C2").Value = "EgyptAir MS-8224"
$ws.Range("D2").Value = 747
$ws.Range("E2").Value = 706
$ws.Range("F2").Value = 41

# --- Row 3 ---
$ws.Range("A3").Value = "'10-JAN-26"
$ws.Range("C3").Value = "EgyptAir MS-8224"
$ws.Range("D3").Value = 747
$ws.Range("F3").Value = 161

# --- Row 4 ---
$ws.Range("A4").Value = "'10-JAN-26"
$ws.Range("B4").Value = "SM-492"
$ws.Range("C4").Value = "EgyptAir MS-8224"
$ws.Range("D4").Value = 747
$ws.Range("E4").Value = 646
$ws.Range("F4").Value = 101

# --- Row 5 ---
$ws.Range("A5").Value = "'16-JAN-26"
$ws.Range("B5").Value = "SM-492"
$ws.Range("C5").Value = "Saudia SV-335"
$ws.Range("D5").Value = 461
$ws.Range("E5").Value = 436
$ws.Range("F5").Value = 25

# --- Row 6 ---
$ws.Range("A6").Value = "'26-JAN-26"
$ws.Range("B6").Value = "SM-492"
$ws.Range("D6").Value = 476
$ws.Range("E6").Value = 471
$ws.Range("F6").Value = 5

# --- New rows 7 & 8: clone formatting from row 6 first ---
$ws.Range("A6:K6").Copy()
$ws.Range("A7:K7").PasteSpecial(-4122)
$ws.Range("A6:K6").Copy()
$ws.Range("A8:K8").PasteSpecial(-4122)

# --- Row 7 (new) ---
$ws.Range("A7").Value = "'29-JAN-26"
$ws.Range("B7").Value = "SM-456"
$ws.Range("C7").Value = "EgyptAir MS-666"
$ws.Range("D7").Value = 586
$ws.Range("E7").Value = 586
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 46
$ws.Range("H7").Value = 30
$ws.Range("I7").Value = -16
$ws.Range("J7").Value = "LOW THREAT"
$ws.Range("K7").Value = "SAR"

# --- Row 8 (new) ---
$ws.Range("A8").Value = "'29-JAN-26"
$ws.Range("B8").Value = "SM-456"
$ws.Range("C8").Value = "EgyptAir MS-662"
$ws.Range("D8").Value = 586
$ws.Range("E8").Value = 586
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 46
$ws.Range("H8").Value = 30
$ws.Range("I8").Value = -16
$ws.Range("J8").Value = "LOW THREAT"
$ws.Range("K8").Value = "SAR"
